# Commit message: "remove Gamelogic project, modify SLG building config"
# Update the SLG building/state function matrix on Sheet1: most flag cells in
# columns C..N for rows 2-13 are cleared to 0, leaving only column B ("always
# available") and column O ("finish") set to 1 (row 2 keeps D=1 as well, per
# the source data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final state of columns B..O for each data row (2..13), taken from the diff.
$rowValues = @{
    2  = @(1,1,1,0,0,0,0,0,0,0,0,0,0,1)
    3  = @(1,0,0,1,0,0,0,0,0,0,0,0,0,1)
    4  = @(1,0,1,0,1,0,0,0,0,0,0,0,0,1)
    5  = @(1,0,0,1,0,0,0,0,0,0,0,0,0,1)
    6  = @(1,0,0,0,0,0,0,0,0,0,0,0,0,1)
    7  = @(1,0,0,0,0,0,0,0,0,0,0,0,0,1)
    8  = @(1,0,0,0,0,0,0,0,0,0,0,0,0,1)
    9  = @(1,0,0,0,0,0,0,0,0,0,0,0,0,1)
    10 = @(1,0,0,0,0,0,0,0,0,0,0,0,0,1)
    11 = @(1,0,0,0,0,0,0,0,0,0,0,0,0,1)
    12 = @(1,0,0,0,0,0,0,0,0,0,0,0,0,1)
    13 = @(1,0,0,0,0,0,0,0,0,0,0,0,0,1)
}

foreach ($r in $rowValues.Keys) {
    $values = $rowValues[$r]
    # Columns B (2) through O (15)
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 2 + $i
        $ws.Cells.Item($r, $col).Value = $values[$i]
    }
}

# Reflect the new active-cell selection recorded in the saved workbook.
$ws.Range("F11").Select()
